# Mark additional attendance cells as "A" (Absent) on the attendance sheet,
# matching the style already used for existing "A" marks (border + centered
# font, readingOrder alignment applied -> style index 36 in the original file).
#
# Cells to be newly marked "A" for each row (column letter -> row number):
#   Row 8 : Q, R
#   Row 12: S
#   Row 18: Q
#   Row 21: S
#   Row 25: Q
#   Row 28: R
#   Row 29: S
#   Row 32: S
#   Row 33: S
#   Row 36: Q, R
#   Row 45: S
#   Row 46: Q
#   Row 47: S
#   Row 48: R, S
#   Row 50: Q
#   Row 51: R
#   Row 54: S
#   Row 59: Q, S
#   Row 60: Q
#   Row 64: S
#   Row 67: S
#   Row 68: Q, R, S
#   Row 69: R
#   Row 72: S

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @(
    "Q8","R8",
    "S12",
    "Q18",
    "S21",
    "Q25",
    "R28",
    "S29",
    "S32",
    "S33",
    "Q36","R36",
    "S45",
    "Q46",
    "S47",
    "R48","S48",
    "Q50",
    "R51",
    "S54",
    "Q59","S59",
    "Q60",
    "S64",
    "S67",
    "Q68","R68","S68",
    "R69",
    "S72"
)

# A cell that already carries the "Absent" formatting (border + centered
# font + alignment) so the newly-marked cells pick up the same cell style
# used elsewhere in the sheet for "A" entries.
$template = $ws.Range("P50")

$targetAddress = [string]::Join(",", $targetCells)
$targetRange = $ws.Range($targetAddress)

foreach ($area in $targetRange.Areas) {
    $template.Copy()
    $area.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

foreach ($addr in $targetCells) {
    $ws.Range($addr).Value = "A"
}
